# Generate Report for Handoff
# The dafdef12-7a35-4b2d-9c12-246e16d27400.md file moved from
# "Handed back: in sync with en-US" to "Ready for handoff", with refreshed
# handoff timestamps and a stale-handback warning recorded for the
# zh-cn / de-de locales.

$wb = $excel.ActiveWorkbook

$newStatus  = "Ready for handoff"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f028c4af5ee1a41b319caa5709aeafef8978a28d/e2e/dafdef12-7a35-4b2d-9c12-246e16d27400.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/23ab20debd7765920c949d5f07ba054c08cd0ef2/e2e/dafdef12-7a35-4b2d-9c12-246e16d27400.md."

# --- Overview sheet: row for dafdef12-...md is row 3 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Range("G3").Value = "2016-08-20 16:56:41"

# --- zh-cn sheet: row for dafdef12-...md is row 3 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("H3").Value = "2016-08-20 16:56:37"
$wsZhCn.Range("P3").Value = $errorDetail

# --- de-de sheet: row for dafdef12-...md is row 3 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("H3").Value = "2016-08-20 16:56:41"
$wsDeDe.Range("P3").Value = $errorDetail
